$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entries finishing section 2 (2.15, 2.16) in rows 20 and 21
$ws.Range("A20").Value = "2.15 Function expressions and arrows"
$ws.Range("A21").Value = "2.16 JavaScript specials"

# Column A grew wider to fit the new longer text (target OOXML width 33.24;
# this engine rounds column widths to whole pixels, so 32.3 is the closest
# achievable ColumnWidth input -> OOXML width 33.166666..).
$ws.Columns.Item(1).ColumnWidth = 32.3

# Selection moved on to the next empty row (19) after finishing row 18
$ws.Rows.Item(19).Select()
